$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New records to append (rows 80-86), matching the target diff.
$data = @(
    @(1119373171, "212417459", "ABDUL SATTAR", "GHULAM AKBAR", "LASHARI", "4120263884305", "B.B.A (HONS)", "ANNUAL FEE", "University of Sindh", 42000, "1BILL_PAY", "2024-01-15T19:00:00.000Z", "12:56:00"),
    @(1119342229, "212417457", "ABDUL HAFEEZ", "ASHIQUE HUSSAIN", "SOOMRO", "4120124622975", "B.B.A (HONS)", "ANNUAL FEE", "University of Sindh", 42000, "1BILL_PAY", "2024-01-15T19:00:00.000Z", "12:19:00"),
    @(1119358164, "212402523", "DUA", "ZULFIQAR ALI", "MEMON", "4160106507616", "B.B.A (HONS)", "ANNUAL FEE", "University of Sindh", 42000, "1BILL_PAY", "2024-01-15T19:00:00.000Z", "12:36:00"),
    @(1119452756, "212414137", "SURAJ KUMAR", "RAJESH KUMAR", "MAHESHWARI", "4440390471307", "BS (COMPUTER SCIENCE) PRE-ENGINEERING", "ANNUAL FEE", "University of Sindh", 42000, "1BILL_PAY", "2024-01-15T19:00:00.000Z", "14:30:00"),
    @(1119483553, "212417462", "ALI RAZA", "GHULAM SARWAR", "PANHWAR", "4120189418207", "B.B.A (HONS)", "ANNUAL FEE", "University of Sindh", 42000, "1BILL_PAY", "2024-01-15T19:00:00.000Z", "15:02:00"),
    @(1119376154, "212414107", "ABDUL MOIZ", "JAHANIGIR AHMED", "SHAIKH", "4130397882001", "BS (COMPUTER SCIENCE) PRE-ENGINEERING", "ANNUAL FEE", "University of Sindh", 42000, "AGENT", "2024-01-15T19:00:00.000Z", "13:01:00"),
    @(1119784016, "212413654", "AMNA KHANUM", "KASHIF ALI KHAN", "YOUSUFZAI PATHAN", "4410327297538", "BS (ARTIFICIAL INTELLIGENCE)", "ANNUAL FEE", "University of Sindh", 42000, "AGENT", "2024-01-15T19:00:00.000Z", "20:08:00")
)

# Columns that must be stored as literal text (even though some look numeric),
# matching t="str" cells in the source sheet. Columns A (1) and J (10) stay numeric.
$textCols = @(2,3,4,5,6,7,8,9,11,12,13)

$startRow = 80
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    # Column A: Tran Id (number)
    $ws.Cells.Item($row, 1).Value = $rowData[0]

    # Text columns: force text number format first so numeric-looking strings
    # (challan numbers, CNIC, times) are preserved as text, not re-parsed as numbers.
    foreach ($col in $textCols) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = [string]$rowData[$col - 1]
    }

    # Column J: Amount (number)
    $ws.Cells.Item($row, 10).Value = $rowData[9]
}
